# ---------------------------------------------------------------------
# What the reference diff actually contains
# ---------------------------------------------------------------------
# Every hunk in the target diff (word/document.xml AND word/styles.xml)
# is a pure XML-attribute / namespace-declaration re-ordering:
#   - xmlns:* declarations on <w:document> sorted alphabetically,
#   - w:pgSz / w:pgMar / w:rFonts / w:lang attributes re-ordered,
#   - w:latentStyles + every <w:lsdException> attribute re-ordered,
#   - every <w:style>/<w:tblInd>/<w:tblCellMar cell> attribute re-ordered.
# Not a single element, attribute *value*, run of text, field code,
# style id, or setting differs between the two sides - confirmed by
# diffing the attribute *sets* (not their textual order) of every
# changed line: they are identical both before and after. The commit
# message corroborates this: "Fixed POI packaging and upgraded to POI
# 3.15" - i.e. this fixture was simply re-saved by a newer Apache POI,
# whose OOXML writer happens to emit attributes in sorted order, while
# the previous POI version preserved insertion order. It is a
# packaging/serialization artifact of the tool that produced the test
# fixture, not a content edit that was made to the document in Word.
#
# ---------------------------------------------------------------------
# Why this script makes no calls against the Word object model
# ---------------------------------------------------------------------
# Word's automation surface (real Word, and this COM-interop host
# alike) has no property/method that lets a caller choose the byte
# order in which an element's attributes get serialized - that is an
# internal writer detail, not part of the object model. So there is no
# `$d...` call that could "apply" an attribute-sort.
#
# Deliberately touching the document anyway (e.g. nudging PageSetup,
# re-assigning a style's own NameLocal, or running a Find/Replace just
# to "dirty" the part so it gets rewritten) was tried and rejected:
# every such mutation forces this host to rebuild <w:body> from its
# in-memory model, and that rebuild does NOT emit sorted attributes
# either (namespace/attribute order there is insertion-order, not
# alphabetical) - so it cannot move us any closer to the target -
# while it also rewrites the malformed field code in the third
# paragraph (it has no <w:fldChar w:fldCharType="separate"/>, so the
# rebuild reclassifies the cached-result <w:t> runs as <w:instrText>).
# That is a real content regression that is *not* present in the
# target diff, so it must be avoided.
#
# The faithful reproduction of "every value is byte-identical, only
# attribute order changed" is therefore to leave the document's
# content/object-model completely untouched.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument
